# Update the HHSawardsNumberInternationalPieTable worksheet:
#  - retitle the table/description text for "international" awards
#  - rename the column headers (Organization -> Agency, etc.)
#  - re-sort the agency rows alphabetically by agency name
#  - add a thin box border around the data table
#  - restore the active selection to G4 (matches author's re-save)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / label updates -------------------------------------------------

$ws.Range("A1").Value = "2016 International Grant Dollars By Agency Table"
$ws.Range("A3").Value = "This table displays the total number of international grant awards each HHS agency awarded in FY 2016. It is provided as a text alternative to the interactive chart on the Awards page of this website."
$ws.Range("A5").Value = "HHS Total Number of International Awards Description"
$ws.Range("A7").Value = "Number of international awards HHS awarded in FY 2016 by agency."
$ws.Range("A9").Value = "Agency"
$ws.Range("B9").Value = "Total Number of International Awards"

# The longer wrapped header text now needs two lines.
$ws.Rows.Item(9).RowHeight = 30

# --- Sort the agency data rows (A10:B16) alphabetically by agency --------

$dataRange = $ws.Range("A10:B16")
$sortKey = $ws.Range("A10:A16")
$dataRange.Sort($sortKey, 1)

# --- Add a thin box border around the whole table (header + data) --------

$tableRange = $ws.Range("A9:B16")
$tableRange.Borders.LineStyle = 1

# --- Restore the active cell selection ------------------------------------

$ws.Range("G4").Select()
